$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.05089966666666667
$ws.Range("H2").Value = 0.152699
$ws.Range("I2").Value = 0.02671091810242436
$ws.Range("J2").Value = 0.03728162213961778
$ws.Range("M2").Value = 12.431794
$ws.Range("N2").Value = 37.295382
$ws.Range("O2").Value = 0.6267040910788743
$ws.Range("P2").Value = 0.7075740515758999
$ws.Range("Q2").Value = 0.6327741706686668
$ws.Range("R2").Value = 5.694967536018001
$ws.Range("S2").Value = 0.01673984165126211
$ws.Range("T2").Value = 0.02637950842665113
$ws.Range("G3").Value = 0.05089966666666667
$ws.Range("H3").Value = 0.152699
$ws.Range("I3").Value = 0.02671091810242436
$ws.Range("J3").Value = 0.03728162213961778
$ws.Range("O3").Value = 0.0264162940991436
$ws.Range("P3").Value = 0.0298250554119953
$ws.Range("Q3").Value = 0.02667215489522222
$ws.Range("R3").Value = 0.240049394057
$ws.Range("S3").Value = 0.0007056034682517807
$ws.Range("T3").Value = 0.001111926446163171
$ws.Range("G4").Value = 0.05089966666666667
$ws.Range("H4").Value = 0.152699
$ws.Range("I4").Value = 0.02671091810242436
$ws.Range("J4").Value = 0.03728162213961778
$ws.Range("M4").Value = 0.03915333333333333
$ws.Range("N4").Value = 0.11746
$ws.Range("O4").Value = 0.001973774193762771
$ws.Range("P4").Value = 0.002228470219130754
$ws.Range("Q4").Value = 0.001992891615555555
$ws.Range("R4").Value = 0.01793602454
$ws.Range("S4").Value = 0.00005272132084227605
$ws.Range("T4").Value = 0.00008308098465902402
$ws.Range("G5").Value = 0.05089966666666667
$ws.Range("H5").Value = 0.152699
$ws.Range("I5").Value = 0.02671091810242436
$ws.Range("J5").Value = 0.03728162213961778
$ws.Range("M5").Value = 6.8015495
$ws.Range("N5").Value = 13.603099
$ws.Range("O5").Value = 0.3428756056708687
$ws.Range("P5").Value = 0.2580802061075034
$ws.Range("Q5").Value = 0.3461966023668334
$ws.Range("R5").Value = 2.077179614201
$ws.Range("S5").Value = 0.009158522222393724
$ws.Range("T5").Value = 0.009621648725814618
$ws.Range("G6").Value = 0.05089966666666667
$ws.Range("H6").Value = 0.152699
$ws.Range("I6").Value = 0.02671091810242436
$ws.Range("J6").Value = 0.03728162213961778
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.04027333333333333
$ws.Range("N6").Value = 0.12082
$ws.Range("O6").Value = 0.002030234957350741
$ws.Range("P6").Value = 0.002292216685470609
$ws.Range("Q6").Value = 0.002049899242222222
$ws.Range("R6").Value = 0.01844909318
$ws.Range("S6").Value = 0.00005422943967447466
$ws.Range("T6").Value = 0.00008545755632984236
$ws.Range("I7").Value = 0.1226793098007496
$ws.Range("J7").Value = 0.1712289953794413
$ws.Range("M7").Value = 12.431794
$ws.Range("N7").Value = 37.295382
$ws.Range("O7").Value = 0.6267040910788743
$ws.Range("P7").Value = 0.7075740515758999
$ws.Range("Q7").Value = 2.906238498418667
$ws.Range("R7").Value = 26.15614648576801
$ws.Range("S7").Value = 0.07688362534286242
$ws.Range("T7").Value = 0.1211571940079023
$ws.Range("I8").Value = 0.1226793098007496
$ws.Range("J8").Value = 0.1712289953794413
$ws.Range("O8").Value = 0.0264162940991436
$ws.Range("P8").Value = 0.0298250554119953
$ws.Range("S8").Value = 0.003240732727576551
$ws.Range("T8").Value = 0.005106914275332124
$ws.Range("I9").Value = 0.1226793098007496
$ws.Range("J9").Value = 0.1712289953794413
$ws.Range("M9").Value = 0.03915333333333333
$ws.Range("N9").Value = 0.11746
$ws.Range("O9").Value = 0.001973774193762771
$ws.Range("P9").Value = 0.002228470219130754
$ws.Range("Q9").Value = 0.009153057448888889
$ws.Range("R9").Value = 0.08237751704
$ws.Range("S9").Value = 0.0002421412557933477
$ws.Range("T9").Value = 0.0003815787168547624
$ws.Range("I10").Value = 0.1226793098007496
$ws.Range("J10").Value = 0.1712289953794413
$ws.Range("M10").Value = 6.8015495
$ws.Range("N10").Value = 13.603099
$ws.Range("O10").Value = 0.3428756056708687
$ws.Range("P10").Value = 0.2580802061075034
$ws.Range("Q10").Value = 1.590029967179333
$ws.Range("R10").Value = 9.540179803076001
$ws.Range("S10").Value = 0.04206374265121616
$ws.Range("T10").Value = 0.04419081441910695
$ws.Range("I11").Value = 0.1226793098007496
$ws.Range("J11").Value = 0.1712289953794413
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.04027333333333333
$ws.Range("N11").Value = 0.12082
$ws.Range("O11").Value = 0.002030234957350741
$ws.Range("P11").Value = 0.002292216685470609
$ws.Range("Q11").Value = 0.009414885075555556
$ws.Range("R11").Value = 0.08473396568000001
$ws.Range("S11").Value = 0.0002490678233011432
$ws.Range("T11").Value = 0.0003924939602451252
$ws.Range("G12").Value = 1.620901
$ws.Range("H12").Value = 3.241802
$ws.Range("I12").Value = 0.8506097720968261
$ws.Range("J12").Value = 0.791489382480941
$ws.Range("M12").Value = 12.431794
$ws.Range("N12").Value = 37.295382
$ws.Range("O12").Value = 0.6267040910788743
$ws.Range("P12").Value = 0.7075740515758999
$ws.Range("Q12").Value = 20.150707326394
$ws.Range("R12").Value = 120.904243958364
$ws.Range("S12").Value = 0.5330806240847498
$ws.Range("T12").Value = 0.5600373491413465
$ws.Range("G13").Value = 1.620901
$ws.Range("H13").Value = 3.241802
$ws.Range("I13").Value = 0.8506097720968261
$ws.Range("J13").Value = 0.791489382480941
$ws.Range("O13").Value = 0.0264162940991436
$ws.Range("P13").Value = 0.0298250554119953
$ws.Range("Q13").Value = 0.8493753569143332
$ws.Range("R13").Value = 5.096252141486
$ws.Range("S13").Value = 0.02246995790331527
$ws.Range("T13").Value = 0.02360621469050001
$ws.Range("G14").Value = 1.620901
$ws.Range("H14").Value = 3.241802
$ws.Range("I14").Value = 0.8506097720968261
$ws.Range("J14").Value = 0.791489382480941
$ws.Range("M14").Value = 0.03915333333333333
$ws.Range("N14").Value = 0.11746
$ws.Range("O14").Value = 0.001973774193762771
$ws.Range("P14").Value = 0.002228470219130754
$ws.Range("Q14").Value = 0.06346367715333333
$ws.Range("R14").Value = 0.3807820629199999
$ws.Range("S14").Value = 0.001678911617127147
$ws.Range("T14").Value = 0.001763810517616968
$ws.Range("G15").Value = 1.620901
$ws.Range("H15").Value = 3.241802
$ws.Range("I15").Value = 0.8506097720968261
$ws.Range("J15").Value = 0.791489382480941
$ws.Range("M15").Value = 6.8015495
$ws.Range("N15").Value = 13.603099
$ws.Range("O15").Value = 0.3428756056708687
$ws.Range("P15").Value = 0.2580802061075034
$ws.Range("Q15").Value = 11.0246383860995
$ws.Range("R15").Value = 44.098553544398
$ws.Range("S15").Value = 0.2916533407972588
$ws.Range("T15").Value = 0.2042677429625818
$ws.Range("G16").Value = 1.620901
$ws.Range("H16").Value = 3.241802
$ws.Range("I16").Value = 0.8506097720968261
$ws.Range("J16").Value = 0.791489382480941
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.04027333333333333
$ws.Range("N16").Value = 0.12082
$ws.Range("O16").Value = 0.002030234957350741
$ws.Range("P16").Value = 0.002292216685470609
$ws.Range("Q16").Value = 0.06527908627333333
$ws.Range("R16").Value = 0.39167451764
$ws.Range("S16").Value = 0.001726937694375123
$ws.Range("T16").Value = 0.001814265168895642
